# Add a new "Sound Cooltime" (cooltime / float) column to both the
# SOUND_BUNDLE and SOUND_RESOURCE tables.
#
# SOUND_BUNDLE (sheet1): new column F
# SOUND_RESOURCE (sheet2): new column G

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# SOUND_BUNDLE sheet - add column F ("cooltime" / "Sound Cooltime")
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Header rows (key row first, then display row, matching the order
# the strings were originally authored in, so the shared-string
# table ends up with "cooltime" before "Sound Cooltime").
$ws1.Range("F2").Value = "cooltime"
$ws1.Range("F1").Value = "Sound Cooltime"
$ws1.Range("F3").Value = "float"

# Data rows.
$ws1.Range("F5").Value = 0.05
$ws1.Range("F6").Value = 0.05
$ws1.Range("F7").Value = 0.05

# Column width + text number format applied last so the numeric
# data values above stay numeric instead of being coerced to text.
$col1 = $ws1.Columns.Item(6)
$col1.ColumnWidth = 8.285714285714286
$col1.NumberFormat = "@"

# ---------------------------------------------------------------
# SOUND_RESOURCE sheet - add column G ("cooltime" / "Sound Cooltime")
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

$ws2.Range("G2").Value = "cooltime"
$ws2.Range("G1").Value = "Sound Cooltime"
$ws2.Range("G3").Value = "float"

$ws2.Range("G5").Value = 0.05
$ws2.Range("G6").Value = 0.05

$col2 = $ws2.Columns.Item(7)
$col2.ColumnWidth = 8.285714285714286
$col2.NumberFormat = "@"
